$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number must keep their original
# Text cell type (matching the source data, which stores prices/labels as
# text, e.g. "211.68"), so force Text number format before assigning those.
$ws.Range("D2").Value = "28.469.68"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "1.564.97"
$ws.Range("E3").Value = "  -2.25%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.68"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.493"
$ws.Range("E6").Value = "  -1.36%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.16"
$ws.Range("E8").Value = "  +4.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.99"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("E10").Value = "  -1.95%  "
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "1.788.63"
$ws.Range("E13").Value = "  -2.24%  "
$ws.Range("D14").Value = "1.553.79"
$ws.Range("E14").Value = "  -2.97%  "
$ws.Range("E15").Value = "  -2.81%  "
$ws.Range("D16").Value = "28.476.74"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("E17").Value = "  -3.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.24"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.45"
$ws.Range("E19").Value = "  -1.99%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.35"
$ws.Range("E20").Value = "  -2.82%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0693"
$ws.Range("E21").Value = "  -2.73%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.89"
$ws.Range("E23").Value = "  -6.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.11"
$ws.Range("E24").Value = "  -3.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.08"
$ws.Range("E25").Value = "  +5.94%  "
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.98"
$ws.Range("E27").Value = "  -2.31%  "
$ws.Range("E28").Value = "  -2.99%  "
$ws.Range("E29").Value = "  -4.20%  "
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("E31").Value = "  -2.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.10"
$ws.Range("E32").Value = "  -3.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.21"
$ws.Range("E33").Value = "  -1.46%  "
$ws.Range("E34").Value = "  -3.08%  "
$ws.Range("D35").Value = "1.390.40"
$ws.Range("E35").Value = "  -2.37%  "
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("E37").Value = "  -3.78%  "
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("E39").Value = "  +2.14%  "
$ws.Range("E40").Value = "  -1.25%  "
$ws.Range("E41").Value = "  -2.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.89"
$ws.Range("E43").Value = "  +1.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.787"
$ws.Range("E44").Value = "  -4.44%  "
$ws.Range("E45").Value = "  -4.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.977"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "62.71"
$ws.Range("E47").Value = "  -3.69%  "
$ws.Range("D48").Value = "1.701.28"
$ws.Range("E48").Value = "  -2.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.95"
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("E50").Value = "  -4.38%  "
$ws.Range("E51").Value = "  -0.60%  "
